$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (rows 9-15), matching the batting-innings records
# that were scraped again and appended to the existing table.
$newRows = @(
    @(" Abu Dhabi", " October 30 2020", "Royals won by 7 wickets (with 15 balls remaining)", "Kings XI Punjab", "Rajasthan Royals", "Mandeep Singh ", "0", "1", "0", "0", "0.00"),
    @(" Sharjah", " October 26 2020", "Kings XI won by 8 wickets (with 7 balls remaining)", "Kings XI Punjab", "Kolkata Knight Riders", "Mandeep Singh ", "66", "56", "8", "2", "117.85"),
    @(" Dubai (DSC)", " October 04 2020", "Super Kings won by 10 wickets (with 14 balls remaining)", "Kings XI Punjab", "Chennai Super Kings", "Mandeep Singh ", "27", "16", "0", "2", "168.75"),
    @(" Abu Dhabi", " November 01 2020", "Super Kings won by 9 wickets (with 7 balls remaining)", "Kings XI Punjab", "Chennai Super Kings", "Mandeep Singh ", "14", "15", "1", "0", "93.33"),
    @(" Abu Dhabi", " October 10 2020", "KKR won by 2 runs", "Kings XI Punjab", "Kolkata Knight Riders", "Mandeep Singh ", "0", "1", "0", "0", "0.00"),
    @(" Dubai (DSC)", " October 24 2020", "Kings XI won by 12 runs", "Kings XI Punjab", "Sunrisers Hyderabad", "Mandeep Singh ", "17", "14", "1", "0", "121.42"),
    @(" Dubai (DSC)", " October 08 2020", "Sunrisers won by 69 runs", "Kings XI Punjab", "Sunrisers Hyderabad", "Mandeep Singh ", "6", "6", "0", "0", "100.00")
)

$startRow = 9
$columns = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K")

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 0; $c -lt $columns.Count; $c++) {
        $cell = $ws.Range($columns[$c] + $row)
        # Force text storage so numeric-looking values (e.g. "66", "0.00")
        # stay as text, matching the rest of the sheet (t="str" cells).
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$c]
    }
}

$ws.Range("A1:K15").Select() | Out-Null
